# Apply the updated crypto price/volume snapshot (and the small block of
# rows 19-26 that got reordered) from the "Updated symbol list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) columns hold plain text in the source file
# (e.g. "328.42" / "1.36%"), not real numbers/percentages. Excel normally
# "smart"-converts such strings into numeric/percentage values when you set
# .Value, which would change both the stored type and the on-disk text. To
# keep these as literal text (matching the original), prepend a leading
# apostrophe, just like typing '328.42 directly into a cell in the UI.
function Set-TextValue($cell, $text) {
    $ws.Range($cell).Value = "'" + $text
}

Set-TextValue "D2" '328.42'
Set-TextValue "E2" '1.36%'
Set-TextValue "D3" '41.64'
Set-TextValue "E3" '5.13%'
Set-TextValue "D4" '5.619'
Set-TextValue "E4" '-4.08%'
Set-TextValue "D5" '0.08172'
Set-TextValue "E5" '1.69%'
Set-TextValue "D6" '2.018'
Set-TextValue "E6" '-1.26%'
Set-TextValue "D7" '8.731'
Set-TextValue "E7" '1.02%'
Set-TextValue "D8" '4.531'
Set-TextValue "E8" '-1.04%'
Set-TextValue "D10" '0.9212'
Set-TextValue "E10" '-1.05%'
Set-TextValue "E11" '0.64%'
Set-TextValue "D12" '0.1953'
Set-TextValue "E12" '0.03%'
Set-TextValue "D13" '0.09414'
Set-TextValue "E13" '2.89%'
Set-TextValue "D14" '0.03803'
Set-TextValue "E14" '6.06%'
Set-TextValue "D15" '0.1058'
Set-TextValue "E15" '1.14%'
Set-TextValue "E16" '0.45%'
Set-TextValue "D17" '0.006279'
Set-TextValue "E17" '1.66%'
$ws.Range("B19").Value = 'LEO'
$ws.Range("C19").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue "D19" '3.441'
Set-TextValue "E19" '2.78%'
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue "D20" '0.3496'
Set-TextValue "E20" '-1.15%'
$ws.Range("B21").Value = 'MCDex'
$ws.Range("C21").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue "D21" '8.290'
Set-TextValue "E21" '-4.68%'
$ws.Range("B22").Value = 'ProBitToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextValue "D22" '0.1393'
Set-TextValue "E22" '1.54%'
$ws.Range("B23").Value = 'ZBToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
Set-TextValue "D23" '0.2413'
Set-TextValue "E23" '-1.44%'
$ws.Range("B24").Value = 'CoinExToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue "D24" '0.04411'
Set-TextValue "E24" '0.12%'
$ws.Range("B25").Value = 'BitKan'
$ws.Range("C25").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextValue "D25" '0.001261'
Set-TextValue "E25" '-0.06%'
$ws.Range("B26").Value = 'HotbitToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextValue "D26" '0.004343'
Set-TextValue "E26" '-1.02%'
Set-TextValue "D27" '0.0001182'
Set-TextValue "E27" '2.79%'
Set-TextValue "D39" '0.02774'
Set-TextValue "E39" '9.17%'
Set-TextValue "D40" '0.05425'
Set-TextValue "E40" '3.61%'
Set-TextValue "D41" '0.007674'
Set-TextValue "E41" '2.88%'
Set-TextValue "D42" '0.1419'
Set-TextValue "E42" '0.94%'
Set-TextValue "D43" '0.009011'
Set-TextValue "E43" '-6.13%'
Set-TextValue "D44" '0.002133'
Set-TextValue "E44" '0.82%'
Set-TextValue "D45" '0.01165'
Set-TextValue "E45" '17.28%'
Set-TextValue "D46" '0.00006643'
Set-TextValue "E46" '-1.50%'
Set-TextValue "E47" '0.17%'
$ws.Range("B48").Value = 'BOLO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
Set-TextValue "D48" '0.003217'
Set-TextValue "E48" '7.24%'
$ws.Range("B49").Value = 'CoinbaseStockToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
Set-TextValue "D49" '0.002283'
Set-TextValue "E49" '-0.31%'
Set-TextValue "E50" '0.17%'
Set-TextValue "D51" '0.0002003'
Set-TextValue "E51" '0.17%'
